# Actualización automática 2025-07-18 11:10:09
#
# Updates sales figures for HIDALGO HIDALGO PEDRO GUSTAVO across the three
# sheets of the workbook, plus the small column-width tweaks that come
# along with Excel's autosave/recalc of the report.

$wb = $excel.ActiveWorkbook

# Excel pads stored column <col width="..."> by exactly 0.8333333333333334
# character-widths relative to the `ColumnWidth` property value, so to make
# the saved width land on an exact integer we subtract that pad back out.
$colPad = 0.8333333333333334

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("H10").Value = 782.1
$ws1.Range("M10").Value = 5283.42
$ws1.Range("M16").Value = 12316.05
$ws1.Range("H22").Value = "3 de 20"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F10").Value = 6304.92
$ws2.Range("F16").Value = 12316.05
$ws2.Range("F22").Value = 36339.92

# Column F width: 13 -> 14
$ws2.Columns.Item(6).ColumnWidth = 14 - $colPad

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D7").Value = 2979
$ws3.Range("E7").Value = -579
$ws3.Range("F7").Value = 1.24125

$ws3.Range("D16").Value = 27237
$ws3.Range("E16").Value = 17029.24
$ws3.Range("F16").Value = 0.6152996052974005

$ws3.Range("D19").Value = 36339.92
$ws3.Range("E19").Value = 29038.07762291769
$ws3.Range("F19").Value = 0.5558432702328797

# Column D width: 14 -> 13 ; Column E width: 23 -> 22
$ws3.Columns.Item(4).ColumnWidth = 13 - $colPad
$ws3.Columns.Item(5).ColumnWidth = 22 - $colPad
